# Update train-data.xlsx: new machine learning (moss) plagiarism scores
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (feature 1) and Column B (feature 2) values, row by row.
# Column C is left untouched - only A/B are refreshed with new scores.
$data = @(
    @{ Row = 2;  A = 0.39;  B = 0.311 }
    @{ Row = 3;  A = 0.425; B = 0.423 }
    @{ Row = 4;  A = 0.003; B = 0.117 }
    @{ Row = 5;  A = 0.013; B = 0 }
    @{ Row = 6;  A = 0.003; B = 0 }
    @{ Row = 7;  A = 0.228; B = 0.228 }
    @{ Row = 8;  A = 0.002; B = 0 }
    @{ Row = 9;  A = 0.347; B = 0.326 }
    @{ Row = 10; A = 0.182; B = 0.171 }
    @{ Row = 11; A = 1;     B = 1 }
    @{ Row = 12; A = 0.002; B = 0.039 }
    @{ Row = 13; A = 0.049; B = 0.172 }
    @{ Row = 14; A = 0;     B = 0 }
    @{ Row = 15; A = 0.01;  B = 0.181 }
    @{ Row = 16; A = 0.103; B = 0.042 }
    @{ Row = 17; A = 0.012; B = 0.009 }
    @{ Row = 18; A = 0.03;  B = 0.105 }
)

foreach ($entry in $data) {
    $ws.Cells.Item($entry.Row, 1).Value = $entry.A
    $ws.Cells.Item($entry.Row, 2).Value = $entry.B
}

# Append a brand-new row 19 with its own A/B/C scores.
$ws.Cells.Item(19, 1).Value = 0.017
$ws.Cells.Item(19, 2).Value = 0
$ws.Cells.Item(19, 3).Value = 0.77
